# Adds the "ODI Bowling Extra" worksheet (mirrors the existing "ODI Batting
# Extra" sheet) and populates it with per-match bowling-extras data:
# MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL.

$wb = $excel.ActiveWorkbook

# New sheet goes after the current last tab ("ODI Batting Extra"), i.e. at
# the very end of the workbook, as sheet #5.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Reuse the exact header formatting (bold, bordered, centered) already used
# by the other data sheets instead of re-creating a one-off style.
$headerSource = $wb.Worksheets.Item("ODI Batting Extra")
$headerSource.Range("A1:C1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)

$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "MAIDEN_OVERS"
$ws.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# The source data (match codes, percentages, ...) is text, not numeric -
# force a text format on the body so values like "20.00%" or "4569" are
# stored verbatim instead of being reinterpreted as numbers/percentages.
$ws.Range("A2:C21").NumberFormat = "@"

$rows = @(
  @(2, '4569', '0', '20.00%'),
  @(3, '4570', '2', '20.00%'),
  @(4, '4572', '0', '30.00%'),
  @(5, '4573', $null, $null),
  @(6, '4575', '0', '10.00%'),
  @(7, '4576', '0', $null),
  @(8, '4578', '1', '10.00%'),
  @(9, '4604', '1', '20.00%'),
  @(10, '4610', '0', $null),
  @(11, '4612', '0', '30.00%'),
  @(12, '4617', '1', '40.00%'),
  @(13, '4625', '0', '10.00%'),
  @(14, '4629', '0', $null),
  @(15, '4632', '0', '10.00%'),
  @(16, '4635', $null, $null),
  @(17, '4677', '0', $null),
  @(18, '4681', '0', $null),
  @(19, '4680', '0', $null),
  @(20, '4702', '0', $null),
  @(21, '4703', '0', $null)
)

foreach ($row in $rows) {
  $r = $row[0]
  if ($row[1] -ne $null) { $ws.Cells.Item($r, 1).Value = $row[1] }
  if ($row[2] -ne $null) { $ws.Cells.Item($r, 2).Value = $row[2] }
  if ($row[3] -ne $null) { $ws.Cells.Item($r, 3).Value = $row[3] }
}

$ws.Range("A1").Select()
